$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.602.37"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  -0.09%  "
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'1.842.46"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  -0.18%  "
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = "'  -2.43%  "
$ws.Range('E4').ClearFormats()
$ws.Range('E5').Value = "'  -1.34%  "
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'1.008"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  -2.11%  "
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = "'0.4307"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'  -1.56%  "
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'0.3727"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'  -1.50%  "
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'0.07290"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  -1.03%  "
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'0.8709"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'  -1.01%  "
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'21.31"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'  -0.83%  "
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'1.856.39"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'  +0.41%  "
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = "'6.718"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'  +0.52%  "
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = "'5.385"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'  -1.82%  "
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = "'0.07128"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'  -0.26%  "
$ws.Range('E15').ClearFormats()
$ws.Range('E16').Value = "'  +4.47%  "
$ws.Range('E16').ClearFormats()
$ws.Range('D18').Value = "'0.000008964"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  -0.76%  "
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = "'1.008"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  -2.11%  "
$ws.Range('E19').ClearFormats()
$ws.Range('E20').Value = "'  -0.58%  "
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'27.602.53"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  -0.14%  "
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'5.180"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  -1.94%  "
$ws.Range('E22').ClearFormats()
$ws.Range('E23').Value = "'  -2.38%  "
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = "'2.068.54"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  -0.78%  "
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'1.965"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'  -4.60%  "
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = "'154.34"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'  -2.50%  "
$ws.Range('E26').ClearFormats()
$ws.Range('E27').Value = "'  -0.57%  "
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'2.157"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  +8.70%  "
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'5.308"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'  +0.15%  "
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'117.52"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  +0.07%  "
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = "'0.08900"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  -1.44%  "
$ws.Range('E31').ClearFormats()
$ws.Range('E32').Value = "'  +0.83%  "
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = "'0.7718"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'  +0.48%  "
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'4.512"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  -0.66%  "
$ws.Range('E34').ClearFormats()
$ws.Range('E35').Value = "'  -3.02%  "
$ws.Range('E35').ClearFormats()
$ws.Range('E36').Value = "'  -2.13%  "
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'1.126"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  -1.75%  "
$ws.Range('E37').ClearFormats()
$ws.Range('E38').Value = "'  +0.06%  "
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = "'0.05297"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  +0.80%  "
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'2.886"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  +1.67%  "
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'7.136"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  +4.30%  "
$ws.Range('E41').ClearFormats()
$ws.Range('D42').Value = "'0.1685"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'  +1.26%  "
$ws.Range('E42').ClearFormats()
$ws.Range('D43').Value = "'0.5107"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "'  -0.97%  "
$ws.Range('E43').ClearFormats()
$ws.Range('D44').Value = "'8.746"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +0.62%  "
$ws.Range('E44').ClearFormats()
$ws.Range('D45').Value = "'10.62"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  -0.64%  "
$ws.Range('E45').ClearFormats()
$ws.Range('E46').Value = "'  -2.75%  "
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'0.4736"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  +1.17%  "
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = "'0.06443"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  -2.37%  "
$ws.Range('E48').ClearFormats()
$ws.Range('E49').Value = "'  -2.15%  "
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = "'1.678"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'  -0.92%  "
$ws.Range('E50').ClearFormats()
$ws.Range('E51').Value = "'  -2.24%  "
$ws.Range('E51').ClearFormats()
